$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("oscillazioni smorzate")
$ws3 = $wb.Worksheets.Item("oscillazioni forzate")

# Update the two input values that drive the recalculated statistics
$ws2.Range("C8").Value = 0.01
$ws2.Range("D8").Value = 0.001

# Add an (empty) formatted cell to the right of F2, matching B2's numeric format
$ws2.Range("G2").NumberFormat = $ws2.Range("B2").NumberFormat

# Extend the shared "2*PI()/A10" style formula down through E16:E17
$ws3.Range("E16").Formula = "=2*PI()/A16"
$ws3.Range("E17").Formula = "=2*PI()/A17"

# Set the active sheet to the second tab ("oscillazioni smorzate")
$ws2.Activate()
